$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text representation (it holds
# values like "1.000" / "29.156.41" that Excel would otherwise coerce to
# numbers and silently reformat). Force text, write the values, then restore
# the default "Normal" style so no stray formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.156.41"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.841.84"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "244.34"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "0.6262"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "0.07507"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "0.2940"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "23.33"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("D11").Value = "0.07719"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "1.874.96"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "5.024"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "0.6768"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").Value = "83.14"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "0.000009289"
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("D17").Value = "5.978"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "29.159.41"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "2.127.12"
$ws.Range("E19").Value = "  +1.98%  "
$ws.Range("D20").Value = "230.77"
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "7.198"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "160.64"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "8.569"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "0.1393"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "1.502"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("D31").Value = "4.157"
$ws.Range("E31").Value = "  +2.01%  "
$ws.Range("D32").Value = "0.05565"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("D33").Value = "1.207"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").Value = "0.7515"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").Value = "1.857"
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Value = "2.662"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "2.772"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").Value = "1.229.01"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").Value = "0.01791"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").Value = "6.577"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").Value = "0.9031"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Value = "2.020.42"
$ws.Range("E44").Value = "  +1.77%  "
$ws.Range("D45").Value = "102.28"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "66.39"
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("D47").Value = "0.00000000122"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "0.5101"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").Value = "0.4093"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "9.138"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").Value = "0.05841"
$ws.Range("E51").Value = "  +1.13%  "

$ws.Range("D2:D51").Style = "Normal"
